{"js": "// Replace each \"NNN\u00f7N=\" division-problem prompt in the worksheet table\n// with its new value, per the commit diff. Every source string is unique\n// within the document, so a plain exact-text search (case/whole-string)\n// for each pair is unambiguous.\nconst replacements = [\n  [\"561\u00f72=\", \"501\u00f74=\"],\n  [\"828\u00f72=\", \"616\u00f78=\"],\n  [\"292\u00f75=\", \"687\u00f78=\"],\n  [\"751\u00f77=\", \"622\u00f78=\"],\n  [\"799\u00f72=\", \"543\u00f72=\"],\n  [\"572\u00f74=\", \"306\u00f74=\"],\n  [\"733\u00f79=\", \"629\u00f73=\"],\n  [\"540\u00f72=\", \"974\u00f76=\"],\n  [\"612\u00f79=\", \"316\u00f75=\"],\n  [\"566\u00f76=\", \"653\u00f74=\"],\n  [\"151\u00f72=\", \"610\u00f76=\"],\n  [\"356\u00f73=\", \"362\u00f79=\"],\n  [\"951\u00f78=\", \"796\u00f78=\"],\n  [\"241\u00f75=\", \"272\u00f76=\"],\n  [\"212\u00f77=\", \"205\u00f75=\"],\n  [\"501\u00f76=\", \"615\u00f75=\"],\n  [\"474\u00f73=\", \"353\u00f74=\"],\n  [\"586\u00f73=\", \"757\u00f76=\"],\n  [\"600\u00f79=\", \"971\u00f77=\"],\n  [\"305\u00f76=\", \"410\u00f72=\"],\n  [\"346\u00f73=\", \"866\u00f75=\"],\n  [\"933\u00f74=\", \"794\u00f77=\"],\n  [\"180\u00f72=\", \"163\u00f77=\"],\n  [\"180\u00f76=\", \"831\u00f74=\"],\n  [\"662\u00f79=\", \"541\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"NNN\u00f7N=\" division-problem prompt in the worksheet table\n# with its new value, per the commit diff. Every source string is unique\n# within the document, so an exact-text Find/Replace (whole match, one\n# occurrence) for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"561\u00f72=\", \"501\u00f74=\"),\n    @(\"828\u00f72=\", \"616\u00f78=\"),\n    @(\"292\u00f75=\", \"687\u00f78=\"),\n    @(\"751\u00f77=\", \"622\u00f78=\"),\n    @(\"799\u00f72=\", \"543\u00f72=\"),\n    @(\"572\u00f74=\", \"306\u00f74=\"),\n    @(\"733\u00f79=\", \"629\u00f73=\"),\n    @(\"540\u00f72=\", \"974\u00f76=\"),\n    @(\"612\u00f79=\", \"316\u00f75=\"),\n    @(\"566\u00f76=\", \"653\u00f74=\"),\n    @(\"151\u00f72=\", \"610\u00f76=\"),\n    @(\"356\u00f73=\", \"362\u00f79=\"),\n    @(\"951\u00f78=\", \"796\u00f78=\"),\n    @(\"241\u00f75=\", \"272\u00f76=\"),\n    @(\"212\u00f77=\", \"205\u00f75=\"),\n    @(\"501\u00f76=\", \"615\u00f75=\"),\n    @(\"474\u00f73=\", \"353\u00f74=\"),\n    @(\"586\u00f73=\", \"757\u00f76=\"),\n    @(\"600\u00f79=\", \"971\u00f77=\"),\n    @(\"305\u00f76=\", \"410\u00f72=\"),\n    @(\"346\u00f73=\", \"866\u00f75=\"),\n    @(\"933\u00f74=\", \"794\u00f77=\"),\n    @(\"180\u00f72=\", \"163\u00f77=\"),\n    @(\"180\u00f76=\", \"831\u00f74=\"),\n    @(\"662\u00f79=\", \"541\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # MatchCase:=True, MatchWholeWord:=True so each unique \"NNN\u00f7N=\" token\n    # is matched exactly once; Forward/Wrap ensure the whole document is\n    # scanned and Replace:=wdReplaceAll (2) applies it.\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
